$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row: the first block of columns (A:J) was generated from the
# "FV2210" input file (previously suffixed "_old"); the second block
# (L:U) was generated from the "FV2304" input file (previously suffixed
# "_new"). Column K ("diff") is unaffected.
$oldSuffixCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
foreach ($col in $oldSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Text -replace "_old$", "_FV2210")
}

$newSuffixCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")
foreach ($col in $newSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Text -replace "_new$", "_FV2304")
}

# Freeze the header row.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a proper Excel Table (ListObject), matching
# the header names so the table columns line up with the renamed headers.
$range = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
